{"js": "// Target: the paragraph that reads\n//   \"GRANT INSERT, SELECT ON lrobinson.profile TO 'proj_user'@'localhost';\"\n// needs to become\n//   \"GRANT INSERT, SELECT, UPDATE ON lrobinson.profile TO 'proj_user'@'localhost';\"\n// i.e. an \", UPDATE\" permission is inserted right after \"GRANT INSERT, SELECT\"\n// (and before \" ON ...\"), which is exactly where Word drops its automatic\n// `_GoBack` bookmark (the last-edit-location bookmark Word always stamps at\n// the point of the most recent text change).\n\nconst body = context.document.body;\n\n// Several paragraphs in this document start with \"GRANT INSERT, SELECT ON \",\n// so search on the short, stable prefix and then disambiguate using the\n// full paragraph text (the one that grants on the \"profile\" table).\nconst hits = body.search(\"GRANT INSERT, SELECT\", { matchCase: true });\nhits.load(\"text\");\nawait context.sync();\n\nconst candidates = hits.items.map((range) => {\n  const para = range.paragraphs.getFirst();\n  para.load(\"text\");\n  return { range, para };\n});\nawait context.sync();\n\nconst target = candidates.find((c) => c.para.text.indexOf(\"lrobinson.profile\") !== -1);\nif (!target) {\n  throw new Error(\"Could not locate the 'GRANT INSERT, SELECT ON lrobinson.profile' paragraph\");\n}\n\n// Collapse to the caret right after \"...SELECT\" (before \" ON ...\") and type\n// the new permission there.\nconst insertionPoint = target.range.getRange(\"After\");\ninsertionPoint.insertText(\", UPDATE\", \"Before\");\nawait context.sync();\n\n// Re-locate the caret right after the text we just inserted and drop Word's\n// \"last edit\" bookmark there, same as the diff shows.\nconst afterUpdate = body.search(\"GRANT INSERT, SELECT, UPDATE\", { matchCase: true });\nafterUpdate.load(\"text\");\nawait context.sync();\n\nconst updated = afterUpdate.items\n  .map((range) => ({ range, para: range.paragraphs.getFirst() }));\nupdated.forEach((u) => u.para.load(\"text\"));\nawait context.sync();\n\nconst updatedTarget = updated.find((u) => u.para.text.indexOf(\"lrobinson.profile\") !== -1);\nif (!updatedTarget) {\n  throw new Error(\"Could not re-locate the updated GRANT paragraph\");\n}\n\nconst bookmarkCaret = updatedTarget.range.getRange(\"After\");\nbookmarkCaret.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Target: the paragraph that reads\n#   \"GRANT INSERT, SELECT ON lrobinson.profile TO 'proj_user'@'localhost';\"\n# needs to become\n#   \"GRANT INSERT, SELECT, UPDATE ON lrobinson.profile TO 'proj_user'@'localhost';\"\n# i.e. an \", UPDATE\" permission is inserted right after \"GRANT INSERT, SELECT\"\n# (and before \" ON ...\"), which is exactly where Word drops its automatic\n# `_GoBack` bookmark (the last-edit-location bookmark Word always stamps at\n# the point of the most recent text change).\n\n$d = $word.ActiveDocument\n\n# Several paragraphs in this document start with \"GRANT INSERT, SELECT ON \",\n# so locate the one that grants on the \"profile\" table specifically.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*GRANT INSERT, SELECT ON lrobinson.profile*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not locate the 'GRANT INSERT, SELECT ON lrobinson.profile' paragraph\"\n}\n\n# Scope Find to this paragraph only and land right after \"...SELECT\".\n$rng = $target.Range\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"GRANT INSERT, SELECT\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$found = $rng.Find.Execute()\nif (-not $found) {\n    throw \"Could not find 'GRANT INSERT, SELECT' inside the target paragraph\"\n}\n\n# Collapse to the caret right after \"...SELECT\" (before \" ON ...\") and type\n# the new permission there.\n$rng.Collapse(0)\n$rng.InsertAfter(\", UPDATE\")\n\n# Re-find the caret right after the text we just inserted and drop Word's\n# \"last edit\" bookmark there, same as the diff shows.\n$rng2 = $target.Range\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"GRANT INSERT, SELECT, UPDATE\"\n$rng2.Find.MatchCase = $true\n$rng2.Find.Forward = $true\n$rng2.Find.Wrap = 0\n$found2 = $rng2.Find.Execute()\nif (-not $found2) {\n    throw \"Could not re-locate the updated GRANT clause\"\n}\n$rng2.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $rng2)\n"}
